# Applies the "output generated at 456a3b4" update to 广州-漫展信息.xlsx
# Sheets: 展览 (Exhibition), 演出 (Performance), 本地生活 (Local Life), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 展览 (Exhibition) -- bump "want to go" counters in column F
# ---------------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 140
$wsExhibit.Range("F6").Value = 642
$wsExhibit.Range("F7").Value = 1207
$wsExhibit.Range("F9").Value = 788
$wsExhibit.Range("F10").Value = 681
$wsExhibit.Range("F15").Value = 894
$wsExhibit.Range("F16").Value = 9702
$wsExhibit.Range("F17").Value = 596
$wsExhibit.Range("F20").Value = 41
$wsExhibit.Range("F23").Value = 1745
$wsExhibit.Range("F24").Value = 23
$wsExhibit.Range("F29").Value = 257
$wsExhibit.Range("F32").Value = 58
$wsExhibit.Range("F33").Value = 95
$wsExhibit.Range("F37").Value = 160

# ---------------------------------------------------------------------------
# Sheet: 演出 (Performance)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 119
# Ticket sold out -- numeric "想去人数" price cell becomes the text "已售罄"
$wsShow.Range("G11").Value = "已售罄"
$wsShow.Range("F16").Value = 262

# ---------------------------------------------------------------------------
# Sheet: 本地生活 (Local Life)
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 806

# ---------------------------------------------------------------------------
# Sheet: 全部类型 (All types) -- aggregated view of the above plus its own
# set of row edits (some counters differ slightly from the per-category
# sheets because the combined sheet was refreshed separately).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 806
$wsAll.Range("F5").Value = 140
$wsAll.Range("F10").Value = 642
$wsAll.Range("F11").Value = 1208
$wsAll.Range("F13").Value = 119
$wsAll.Range("F14").Value = 788
$wsAll.Range("F15").Value = 681
$wsAll.Range("F19").Value = 894
$wsAll.Range("F20").Value = 9702
$wsAll.Range("F22").Value = 596
$wsAll.Range("F24").Value = 41
$wsAll.Range("F26").Value = 1745
$wsAll.Range("F27").Value = 23

# Row 29 on "全部类型" used to be the LoveLive anniversary tour; it is
# replaced by the event that used to sit in row 30 (Project SEKAI 25-ji
# doujin tea party x Akiyama Mizuki birthday).
$wsAll.Range("C29").Value = "广州·世界计划25时主题同人茶会×晓山瑞希生日会"
$wsAll.Range("D29").Value = "黄边地铁B出口黄边美食广场1层 胡桃里音乐馆(黄边店)"
$wsAll.Range("E29").Value = "2024.08.24 10:00-08.24 16:30"
$wsAll.Range("F29").Value = 174
$wsAll.Range("G29").Value = 58
$wsAll.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=87815"
$wsAll.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202406/rzS5X2Ko1718735908971.png"

# Row 30 is replaced by the event that used to sit in row 32 (Spring Day
# Project 2024 unplugged concert).
$wsAll.Range("C30").Value = "广州·春日计划2024——特别二次元不插电音乐会"
$wsAll.Range("D30").Value = "人民北路696号 广州友谊剧院"
$wsAll.Range("E30").Value = "2024.08.24 19:30-08.24 21:00"
$wsAll.Range("F30").Value = 80
$wsAll.Range("G30").Value = 88
$wsAll.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=89964"
$wsAll.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202407/lHPV2n6t1722233858047.jpeg"

# Row 32 is replaced entirely by a new event (Meng Jinghui's classic play,
# starring Huang Xiangli, "Letter from an Unknown Woman").
# B32 holds a plain text date string (like the rest of column B) -- force
# text so Excel doesn't reinterpret "2024-08-30" as a real date serial,
# then drop the formatting change so the cell keeps its original (default)
# style, same as every other cell in column B.
$wsAll.Range("B32").NumberFormat = "@"
$wsAll.Range("B32").Value = "2024-08-30"
$wsAll.Range("B32").ClearFormats()
$wsAll.Range("C32").Value = "广州·孟京辉经典戏剧作品·黄湘丽主演《一个陌生女人的来信》"
$wsAll.Range("D32").Value = "广州市越秀区人民北路696号 广州友谊剧院"
$wsAll.Range("E32").Value = "2024.08.30 19:30-08.31 16:30"
$wsAll.Range("F32").Value = 16
$wsAll.Range("G32").Value = 100
$wsAll.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=84570"
$wsAll.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202404/SscDFm1z1713177818070.jpeg"

$wsAll.Range("F36").Value = 257
$wsAll.Range("F39").Value = 58
$wsAll.Range("F40").Value = 95
$wsAll.Range("F47").Value = 160
